$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rngD = $ws.Range("D2:D51")
$rngD.NumberFormat = "@"

$ws.Range("D2").Value = "26.565.89"
$ws.Range("E2").Value = "  +0.18%  "
$ws.Range("D3").Value = "1.627.99"
$ws.Range("E3").Value = "  +0.87%  "
$ws.Range("E4").Value = "  -0.40%  "
$ws.Range("E5").Value = "  -0.17%  "
$ws.Range("D6").Value = "0.502"
$ws.Range("E6").Value = "  -0.55%  "
$ws.Range("E7").Value = "  -0.36%  "
$ws.Range("E8").Value = "  +0.36%  "
$ws.Range("E9").Value = "  +0.05%  "
$ws.Range("D10").Value = "19.24"
$ws.Range("E10").Value = "  +0.28%  "
$ws.Range("D11").Value = "0.0857"
$ws.Range("E11").Value = "  -0.21%  "
$ws.Range("D12").Value = "1.857.43"
$ws.Range("E12").Value = "  +0.66%  "
$ws.Range("D13").Value = "1.633.50"
$ws.Range("E13").Value = "  +1.01%  "
$ws.Range("D14").Value = "4.05"
$ws.Range("E14").Value = "  +0.61%  "
$ws.Range("E15").Value = "  +0.98%  "
$ws.Range("D16").Value = "63.94"
$ws.Range("E16").Value = "  -0.58%  "
$ws.Range("D17").Value = "26.579.89"
$ws.Range("E17").Value = "  +0.10%  "
$ws.Range("D18").Value = "234.02"
$ws.Range("E18").Value = "  +1.21%  "
$ws.Range("D19").Value = "7.81"
$ws.Range("E19").Value = "  +2.13%  "
$ws.Range("D20").Value = "0.0₃0725"
$ws.Range("E20").Value = "  -0.47%  "
$ws.Range("D21").Value = "0.999"
$ws.Range("E21").Value = "  -0.37%  "
$ws.Range("E22").Value = "  -0.42%  "
$ws.Range("D23").Value = "9.17"
$ws.Range("E23").Value = "  +1.57%  "
$ws.Range("E24").Value = "  -0.01%  "
$ws.Range("D25").Value = "145.98"
$ws.Range("E25").Value = "  +0.18%  "
$ws.Range("E26").Value = "  -0.42%  "
$ws.Range("E27").Value = "  +1.09%  "
$ws.Range("E28").Value = "  -0.11%  "
$ws.Range("E29").Value = "  +1.02%  "
$ws.Range("D30").Value = "0.0495"
$ws.Range("E30").Value = "  -0.04%  "
$ws.Range("E31").Value = "  -0.44%  "
$ws.Range("D32").Value = "1.526.30"
$ws.Range("E32").Value = "  +4.86%  "
$ws.Range("E33").Value = "  +0.66%  "
$ws.Range("D34").Value = "3.02"
$ws.Range("E34").Value = "  +1.21%  "
$ws.Range("D35").Value = "1.53"
$ws.Range("E35").Value = "  +3.40%  "
$ws.Range("D36").Value = "2.42"
$ws.Range("E36").Value = "  -0.49%  "
$ws.Range("E37").Value = "  +0.62%  "
$ws.Range("E38").Value = "  -0.15%  "
$ws.Range("D39").Value = "0.837"
$ws.Range("E39").Value = "  +0.33%  "
$ws.Range("D40").Value = "5.87"
$ws.Range("E40").Value = "  -0.55%  "
$ws.Range("E41").Value = "  -0.31%  "
$ws.Range("E42").Value = "  -2.96%  "
$ws.Range("D43").Value = "1.769.11"
$ws.Range("E43").Value = "  +0.78%  "
$ws.Range("D44").Value = "63.24"
$ws.Range("E44").Value = "  +2.89%  "
$ws.Range("E45").Value = "  -0.35%  "
$ws.Range("E46").Value = "  -3.82%  "
$ws.Range("D47").Value = "89.68"
$ws.Range("E47").Value = "  +1.79%  "
$ws.Range("E48").Value = "  +1.67%  "
$ws.Range("D49").Value = "0.0⁦0102"
$ws.Range("E49").Value = "  +1.82%  "
$ws.Range("E50").Value = "  -0.21%  "
$ws.Range("D51").Value = "0.0966"
$ws.Range("E51").Value = "  +0.90%  "

$rngD.Style = "Normal"
